$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expenditure")
foreach ($pt in $ws.PivotTables()) {
    try {
        $pt.RefreshTable()
        Write-Host "refreshed" $pt.Name
    } catch {
        Write-Host "ERROR" $_
    }
}
